# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Price (col D) and Volume(1h) (col E) cells are stored as text in the sheet
# (e.g. "26.277.36", "  -1.67%  "), so every write below targets .Value with a
# string literal. Some new Price values happen to look like plain decimals
# (e.g. "209.77"); Excel would otherwise auto-coerce those to numbers, so for
# those specific cells we force Text format, write the string, then clear the
# number format again (the value/type is what must change, not the formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.277.36'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '1.584.59'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.77'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('E6').Value = '  -1.18%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  -1.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.245'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.59'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0846'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('D12').Value = '1.807.63'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').Value = '1.580.82'
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('E15').Value = '  -1.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.59'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.95%  '
$ws.Range('D17').Value = '26.274.66'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.20'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '206.88'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.94%  '
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.23'
$ws.Range('D23').ClearFormats()
$ws.Range('E24').Value = '  -1.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.61'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.25'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0503'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.96%  '
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.23'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.30'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +11.29%  '
$ws.Range('D35').Value = '1.284.79'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.607'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0167'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.52%  '
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.769'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.80%  '
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.35'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.49%  '
$ws.Range('D45').Value = '1.720.40'
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.89'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.41%  '
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').Value = '0.0₇0996'
$ws.Range('E49').Value = '  -6.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0509'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.70%  '
$ws.Range('E51').Value = '  -0.02%  '
